# The deck originally carries two theme parts:
#   ppt/theme/theme1.xml -> "Integral"     (bound to the slide master / Design)
#   ppt/theme/theme2.xml -> "Office Theme" (bound to the notes master)
# The commit swaps their content, so the slide master ends up using the
# plain "Office Theme" palette and the notes master ends up with the old
# "Integral" palette. The only part of that swap reachable through the
# PowerPoint object model (slide/design theme colours) is applied here via
# ThemeColorScheme - each of the 12 theme colour slots is rewritten from
# the "Integral" values to the "Office" values, in clrScheme document
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
